# Adds an "Exchange Rates" sheet for capital calls upload.
# Moves columns J:M (From Currency, To Currency, Exchange Rate, As Of) off
# the CapitalCall sheet into a new "Exchange Rates" sheet as columns A:D.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("CapitalCall")

# Capture the exchange-rate block (header + 3 data rows) before removing it.
$values = $ws1.Range("J1:M4").Value2

# Add the new worksheet after the existing one and rename it.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Exchange Rates"

# Write header + data into the new sheet.
$ws2.Range("A1:D4").Value2 = $values

# Copy the date number format used on column M (As Of) to column D.
$ws2.Range("D2:D4").NumberFormat = $ws1.Range("M2:M4").NumberFormat

# Remove the now-duplicated columns from the CapitalCall sheet.
$ws1.Range("J1:M4").Clear()

# Update selections to match the target state.
[void]$ws1.Range("C26").Select()
[void]$ws2.Range("D31").Select()

[void]$ws1.Activate()
